$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new columns of data (entered in reading order so the shared-string
# table is built in the same order as the target workbook)
$ws.Range("E1").Value = "Error message"
$ws.Range("E2").Value = "Please complete all required fields before sending."
$ws.Range("F1").Value = "Sign in txt"
$ws.Range("F2").Value = "Thank you"

# Size the new column
$ws.Columns.Item(5).ColumnWidth = 12.8

# Update the active selection
$ws.Range("E2").Select()

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
